# Auto-generated PowerShell-style Excel COM-interop script
# Applies cell-level value updates to the cryptos worksheet per the commit diff.
# Column D (Price) and E (Volume) cells hold text that often *looks* numeric
# (e.g. '1.004', '27.964.44'), so we force NumberFormat='@' (Text) before
# assigning the Value, matching how the source data is stored as inline strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.964.44'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.55'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5111'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3799'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08354'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -8.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.106'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.28'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.192'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.867.49'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.41'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.164'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001089'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.11'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06631'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.78'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.002'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.996.71'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.253'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.563'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.080.06'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.21'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.43'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.81'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1055'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.037'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.567'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.601'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.541'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06506'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02406'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2146'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.207'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6368'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.225'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.25'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.847'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6013'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.89'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.284'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.657'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.966'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.204'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.38'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.48'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.03%  '
